$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting period / validation / update dates on row 8
$ws.Range("B8").Value = 44743
$ws.Range("C8").Value = 44926
$ws.Range("I8").Value = 44936
$ws.Range("J8").Value = 44936
